$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.276082992553711
$ws.Range("B1").Value = 3.283247232437134
$ws.Range("C1").Value = 2.538605928421021
$ws.Range("D1").Value = 1.325958609580994
$ws.Range("E1").Value = 0.9827471375465393
